$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.915.40'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.552.45'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.03'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.39%  '
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0855'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('D12').Value = '1.773.18'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').Value = '1.547.19'
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.520'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').Value = '26.906.85'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('E18').Value = '  +3.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.65'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('E29').Value = '  -0.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0470'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.39%  '
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.11'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.61%  '
$ws.Range('D34').Value = '1.410.49'
$ws.Range('E34').Value = '  +0.54%  '
$ws.Range('E35').Value = '  +2.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.973'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0166'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.525'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.808'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('E42').Value = '  +3.23%  '
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('D47').Value = '1.687.17'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0522'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.83%  '
$ws.Range('D50').Value = '0.0₆0101'
$ws.Range('E50').Value = '  +4.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0960'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.08%  '
